$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 24.38000000000037
$ws.Range("H2").Value = [double]"2.801824667823739e-16"
$ws.Range("K2").Value = 39.61175967413118
$ws.Range("L2").Value = "[31.584860987286632, 47.638658360975725]"
$ws.Range("O2").Value = 1.603816069400194
$ws.Range("P2").Value = "[1.377394977249578, 1.8302371615508104]"
$ws.Range("S2").Value = 68.08541008979232
$ws.Range("T2").Value = "[62.77513529957389, 73.39568488001075]"
$ws.Range("W2").Value = 18.15687687687716
$ws.Range("X2").Value = 17.27831831831858
$ws.Range("Y2").Value = 19.03543543543573

# Row 3 updates
$ws.Range("E3").Value = 25.40000000000053
$ws.Range("H3").Value = [double]"2.801824667823739e-16"
$ws.Range("I3").Value = 0.4614286843722659
$ws.Range("K3").Value = 51.57627642954376
$ws.Range("L3").Value = "[39.62739804868923, 63.52515481039829]"
$ws.Range("M3").Value = [double]"6.661338147750939e-16"
$ws.Range("N3").Value = [double]"6.661338147750939e-16"
$ws.Range("O3").Value = -2.339684618889696
$ws.Range("P3").Value = "[-2.578684660604234, -2.100684577175157]"
$ws.Range("S3").Value = 63.90265929311757
$ws.Range("T3").Value = "[57.68691314622927, 70.11840544000587]"
$ws.Range("W3").Value = 9.458258258258455
$ws.Range("X3").Value = 8.492092092092271
$ws.Range("Y3").Value = 10.42442442442464
